$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 121; this shifts the existing rows 121-164 down
# to 122-165, preserving all of their data (matching the target diff, where
# every row from 121 onward is effectively the previous row's data, and a
# brand-new record appears at row 121 while the last old record lands on the
# new row 165).
$ws.Rows("121:121").Insert()

# Populate the newly inserted row 121 with its data.
$ws.Range("A121").Value = 9
$ws.Range("B121").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C121").Value = "Metropolitana"
$ws.Range("D121").Value = 44466
$ws.Range("E121").Value = 13
$ws.Range("F121").Value = 100112021
$ws.Range("G121").Value = "Ají"
$ws.Range("H121").Value = "Inferno"
$ws.Range("I121").Value = "Primera"
$ws.Range("J121").Value = 25
$ws.Range("K121").Value = 48000
$ws.Range("L121").Value = 50000
$ws.Range("M121").Value = 48960
$ws.Range("N121").Value = "$/caja 12 kilos"
$ws.Range("O121").Value = "Región de Arica y Parinacota"
$ws.Range("P121").Value = 4080
$ws.Range("Q121").Value = 12
$ws.Range("R121").Value = "Hortaliza"
